$d = $word.ActiveDocument
$bm = $d.Bookmarks("_Toc90500056")
$rng = $bm.Range
$rng.Collapse(1)
$insertPos = $rng.Start
$rng.InsertParagraphBefore()
$newp = $d.Range($insertPos, $insertPos + 1)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Rủi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Tên </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rủi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Thiếu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> người </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lập</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trình</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Mô</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tả</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rủi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Dự</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>án</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>th</w:t></w:r><w:r><w:t>iếu</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> người </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>có</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kinh</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>nghiệm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trong</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> các </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>mảng</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>khác</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> nhau</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Xác</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>suất</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>xảy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ra</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cao</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Mức</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>độ</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thiệt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>hại</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cao</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="40"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Giải</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pháp</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> xử lý: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tuy</w:t></w:r><w:r><w:t>ển</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>thêm</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>lập</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>trình</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>viên</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newp.InsertXML($xml)
